$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3197726666666666
$ws.Range("H2").Value = 0.959318
$ws.Range("I2").Value = 0.1074590987069417
$ws.Range("J2").Value = 0.1074590987069417
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 0.01908445911022222
$ws.Range("R2").Value = 0.171760131992
$ws.Range("S2").Value = 0.002796889168061823
$ws.Range("T2").Value = 0.002796889168061823
$ws.Range("G3").Value = 0.3197726666666666
$ws.Range("H3").Value = 0.959318
$ws.Range("I3").Value = 0.1074590987069417
$ws.Range("J3").Value = 0.1074590987069417
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 0.1055960761228889
$ws.Range("R3").Value = 0.9503646851059999
$ws.Range("S3").Value = 0.0154754462671539
$ws.Range("T3").Value = 0.01547544626715389
$ws.Range("G4").Value = 0.3197726666666666
$ws.Range("H4").Value = 0.959318
$ws.Range("I4").Value = 0.1074590987069417
$ws.Range("J4").Value = 0.1074590987069417
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 0.6085622398873333
$ws.Range("R4").Value = 5.477060158986
$ws.Range("S4").Value = 0.08918676327172596
$ws.Range("T4").Value = 0.08918676327172596
$ws.Range("I5").Value = 0.7259488187057992
$ws.Range("J5").Value = 0.7259488187057991
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 0.1289266401208889
$ws.Range("R5").Value = 1.160339761088
$ws.Range("S5").Value = 0.01889461583092884
$ws.Range("T5").Value = 0.01889461583092884
$ws.Range("I6").Value = 0.7259488187057992
$ws.Range("J6").Value = 0.7259488187057991
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("S6").Value = 0.1045456557124438
$ws.Range("T6").Value = 0.1045456557124438
$ws.Range("I7").Value = 0.7259488187057992
$ws.Range("J7").Value = 0.7259488187057991
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("S7").Value = 0.6025085471624265
$ws.Range("T7").Value = 0.6025085471624264
$ws.Range("I8").Value = 0.1665920825872592
$ws.Range("J8").Value = 0.1665920825872591
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05968133333333333
$ws.Range("N8").Value = 0.179044
$ws.Range("O8").Value = 0.02602747651633847
$ws.Range("P8").Value = 0.02602747651633848
$ws.Range("Q8").Value = 0.02958632471777778
$ws.Range("R8").Value = 0.26627692246
$ws.Range("S8").Value = 0.004335971517347808
$ws.Range("T8").Value = 0.004335971517347807
$ws.Range("I9").Value = 0.1665920825872592
$ws.Range("J9").Value = 0.1665920825872591
$ws.Range("O9").Value = 0.144012433133819
$ws.Range("P9").Value = 0.144012433133819
$ws.Range("S9").Value = 0.02399133115422132
$ws.Range("T9").Value = 0.02399133115422131
$ws.Range("I10").Value = 0.1665920825872592
$ws.Range("J10").Value = 0.1665920825872591
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("S10").Value = 0.1382647799156901
$ws.Range("T10").Value = 0.13826477991569
